# "V model working with V based Decision model"
# Adds a new column F of data (values + a "result" column) to the small
# Sheet1 table, renumbers a couple of existing cells, and moves the
# selection to the new last cell (F3). Also restores the column widths
# so every column in the table keeps an explicit (custom) width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers / first data row) ---------------------------------
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
$ws.Range("D1").Value = 4
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 6

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = 35932
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.857

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = 39634
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 491.49

# --- Column widths --------------------------------------------------
# Re-apply an explicit width to every used column (1-14) so each one is
# flagged as a custom width, and give column A / F the wider size used
# for the "label" columns (to match columns A and J).
$ws.Columns.Item(1).ColumnWidth  = 5.833333333333333
$ws.Columns.Item(2).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(3).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(4).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(5).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(6).ColumnWidth  = 6
$ws.Columns.Item(7).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(8).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(9).ColumnWidth  = 1.8333333333333333
$ws.Columns.Item(10).ColumnWidth = 5.833333333333333
$ws.Columns.Item(11).ColumnWidth = 1.8333333333333333
$ws.Columns.Item(12).ColumnWidth = 1.8333333333333333
$ws.Columns.Item(13).ColumnWidth = 1.8333333333333333
$ws.Columns.Item(14).ColumnWidth = 1.8333333333333333

# --- Selection ------------------------------------------------------
# Move the active selection to the new bottom-right cell of the table.
$ws.Range("F3").Select()
